$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit adds a "count rows" test, and in the process the cell that
# previously held a stray aggregate value (22) is corrected to the real
# data-row value (2); the active selection also moves onto that cell.
$ws.Range("B4").Value = 2
$ws.Range("B4").Select()
